# Updates the title-slide credits text box (Shape 2 on Slide 1):
#  - "Marcio Bedran M. da Costa"      -> split into 3 runs ("Marcio " / "Bedran" / " M. da Costa")
#  - "Thiago Rodrigues da Motta Fagundes" -> "Thiago R. da Motta Fagundes"
#  - "Victor Verdan Braga"            -> split into 3 runs ("Victor " / "Verdan" / " Braga")

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Paragraph 10: "Marcio Bedran M. da Costa" -> "Marcio " + "Bedran" + " M. da Costa"
$paraMarcio = $tr.Paragraphs(10, 1)
$midMarcio = $paraMarcio.Characters(8, 6)
$midMarcio.Text = "Bedran"

# Paragraph 12: "Thiago Rodrigues da Motta Fagundes" -> "Thiago R. da Motta Fagundes"
$paraThiago = $tr.Paragraphs(12, 1)
$runThiago = $paraThiago.Runs(1, 1)
$runThiago.Text = "Thiago R. da Motta Fagundes"

# Paragraph 13: "Victor Verdan Braga" -> "Victor " + "Verdan" + " Braga"
$paraVictor = $tr.Paragraphs(13, 1)
$midVictor = $paraVictor.Characters(8, 6)
$midVictor.Text = "Verdan"
